# Refined metadata to be additional tab
$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# --- 1. Update time_taken (column F) timestamps on the "data" sheet ---
$timestamps = @(
    "2021-10-05 14:22:27.314022",
    "2021-10-05 14:22:27.314028",
    "2021-10-05 14:22:27.314031",
    "2021-10-05 14:22:27.314033",
    "2021-10-05 14:22:27.314035",
    "2021-10-05 14:22:27.314037",
    "2021-10-05 14:22:27.314039",
    "2021-10-05 14:22:27.314041",
    "2021-10-05 14:22:27.314043",
    "2021-10-05 14:22:27.314045",
    "2021-10-05 14:22:27.314047",
    "2021-10-05 14:22:27.314049",
    "2021-10-05 14:22:27.314051",
    "2021-10-05 14:22:27.314053",
    "2021-10-05 14:22:27.314055",
    "2021-10-05 14:22:27.314057",
    "2021-10-05 14:22:27.314059",
    "2021-10-05 14:22:27.314061",
    "2021-10-05 14:22:27.314063",
    "2021-10-05 14:22:27.314065",
    "2021-10-05 14:22:27.314067",
    "2021-10-05 14:22:27.314069",
    "2021-10-05 14:22:27.314071"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $dataSheet.Cells.Item($row, 6).Value = $timestamps[$i]
}

# --- 2. Add a new "metadata" worksheet positioned after "data" ---
$metaSheet = $wb.Worksheets.Add($null, $dataSheet)
$metaSheet.Name = "metadata"

# Header row text
$metaSheet.Cells.Item(1, 2).Value = "data_name"
$metaSheet.Cells.Item(1, 3).Value = "data_id"
$metaSheet.Cells.Item(1, 4).Value = "data_version"
$metaSheet.Cells.Item(1, 5).Value = "data_version_created"
$metaSheet.Cells.Item(1, 6).Value = "panel_query_time"
$metaSheet.Cells.Item(1, 7).Value = "panel_get_request"

# Carry over the bold/centered/bordered header style used on the "data" sheet's
# header row (B1:F1, style index referenced as "s=1") onto the metadata header
# row B1:G1, and the same style used on "data"!A2 onto metadata!A2.
$dataSheet.Range("B1:F1").Copy()
$metaSheet.Range("B1:G1").PasteSpecial(-4122)  # xlPasteFormats
$dataSheet.Range("A2").Copy()
$metaSheet.Range("A2").PasteSpecial(-4122)     # xlPasteFormats

# Data row
$metaSheet.Cells.Item(2, 1).Value = 0
$metaSheet.Cells.Item(2, 2).Value = "RASopathies"
$metaSheet.Cells.Item(2, 3).Value = 48

# data_version ("1.75") must be stored as text, not a number. Build it as a
# text formula in a scratch cell, then paste-special the computed value
# (still text) into D2 so no numeric coercion / numberformat style is applied.
$scratch = $metaSheet.Cells.Item(100, 100)
$scratch.Formula = '="1.75"'
$scratch.Copy()
$metaSheet.Cells.Item(2, 4).PasteSpecial(-4163)  # xlPasteValues
$scratch.Clear()

$metaSheet.Cells.Item(2, 5).Value = "2021-01-29T15:21:16.036291Z"
$metaSheet.Cells.Item(2, 6).Value = "2021-10-05 14:22:27.311653"
$metaSheet.Cells.Item(2, 7).Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/48/?format=json"

Write-Output "done"
